$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now filled in
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicate "Contact" / "No display for ContactDetail" row; it becomes
# "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 (the second duplicate "Contact" row) is removed entirely, shifting
# everything below it up by one row.
$meta.Rows.Item(11).Delete()

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# Short column (K) and Definition column (L) for the root Extension row (row 2)
$elements.Range("K2").Value = "Employee Job Grade"
$elements.Range("L2").Value = "Code for the job grade of the employee"
